# Regenerate the localization-status report: two files ("3a4db5f3-..." and
# "db0e1d4f-...") move from "Ready for handoff" to "In Translation" (they
# picked up a fresh handoff round), and two brand-new files
# ("09caa6ba-..." and "656b4403-...") show up as newly "Ready for handoff".
# The ".localization-config" row stays last on every sheet.

$wb = $excel.ActiveWorkbook

$mdSha    = "dba06034e0ad3db0965ae71e6ce5e8b3203defbd"
$zhXlfSha = "8d1b4adfbb98e82488b2259dabf899b0cf9ebc35"
$deXlfSha = "412ada91a8077554672b6e943c17c728057c5ce8"

function MdUrl($name) {
    return "https://github.com/OpenLocalizationTest/oltest/blob/$mdSha/e2e/$name"
}
function ZhXlfUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhXlfSha/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$name"
}
function DeXlfUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deXlfSha/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$name"
}
function XlfUrl($locale, $name) {
    if ($locale -eq "zh") {
        return ZhXlfUrl $name
    } else {
        return DeXlfUrl $name
    }
}
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$mdSha/.localization-config"

$file1 = "3a4db5f3-fc9b-4cba-b03e-702045f93245.md"
$file2 = "db0e1d4f-21a9-457e-be0d-8b9178ecaf07.md"
$file3 = "09caa6ba-7754-4c23-b998-8e4b62501455.md"
$file4 = "656b4403-eba8-4361-af3f-dc3c9dfc61f7.md"

$xlf1zh = "3a4db5f3-fc9b-4cba-b03e-702045f93245.f9416a1d13b0275d82cca981d1eaca99c5bf4127.zh-cn.xlf"
$xlf2zh = "db0e1d4f-21a9-457e-be0d-8b9178ecaf07.08da40945bd1da4b3d9f5e26651018b6748af365.zh-cn.xlf"
$xlf3zh = "09caa6ba-7754-4c23-b998-8e4b62501455.e7133d6e729137647269452931a8ff64cfb5ba8b.zh-cn.xlf"
$xlf4zh = "656b4403-eba8-4361-af3f-dc3c9dfc61f7.4d6233e0638f1dcf647ea351b67ed4f0741b5a07.zh-cn.xlf"

$xlf1de = "3a4db5f3-fc9b-4cba-b03e-702045f93245.f9416a1d13b0275d82cca981d1eaca99c5bf4127.de-de.xlf"
$xlf2de = "db0e1d4f-21a9-457e-be0d-8b9178ecaf07.08da40945bd1da4b3d9f5e26651018b6748af365.de-de.xlf"
$xlf3de = "09caa6ba-7754-4c23-b998-8e4b62501455.e7133d6e729137647269452931a8ff64cfb5ba8b.de-de.xlf"
$xlf4de = "656b4403-eba8-4361-af3f-dc3c9dfc61f7.4d6233e0638f1dcf647ea351b67ed4f0741b5a07.de-de.xlf"

$inTranslation = "In Translation"
$readyHandoff  = "Ready for handoff"
$notLocalized  = "Not to be localized"
$included      = "Include"
$ignored       = "Ignored"
$epoch         = "0001-01-01 00:00:00"
$zhTime1       = "2016-01-18 02:26:05"
$zhTime3       = "2016-01-18 02:27:46"
$deTime1       = "2016-01-18 02:26:19"
$deTime3       = "2016-01-18 02:27:58"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $file1
$ov.Range("B2").Value = $inTranslation
$ov.Range("C2").Value = $inTranslation

$ov.Range("A3").Value = $file2
$ov.Range("B3").Value = $inTranslation
$ov.Range("C3").Value = $inTranslation

$ov.Range("A4").Value = $file3
$ov.Range("B4").Value = $readyHandoff
$ov.Range("C4").Value = $readyHandoff

$ov.Range("A5").Value = $file4
$ov.Range("B5").Value = $readyHandoff
$ov.Range("C5").Value = $readyHandoff

$ov.Range("A6").Value = ".localization-config"
$ov.Range("B6").Value = $notLocalized
$ov.Range("C6").Value = $notLocalized

$ov.Range("A2:A6").Style = "Hyperlink"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), (MdUrl $file1), $null, $null, $file1) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), (MdUrl $file2), $null, $null, $file2) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), (MdUrl $file3), $null, $null, $file3) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A5"), (MdUrl $file4), $null, $null, $file4) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A6"), $configUrl, $null, $null, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (and identically shaped "de-de")
# ---------------------------------------------------------------------------
function FillDetailSheet($ws, $xlf1, $xlf2, $xlf3, $xlf4, $t1, $t3, $locale) {

    $ws.Range("A2").Value = $file1
    $ws.Range("B2").Value = $inTranslation
    $ws.Range("C2").Value = $xlf1
    $ws.Range("D2").Value = $t1
    $ws.Range("G2").Value = $epoch
    $ws.Range("H2").Value = $included

    $ws.Range("A3").Value = $file2
    $ws.Range("B3").Value = $inTranslation
    $ws.Range("C3").Value = $xlf2
    $ws.Range("D3").Value = $t1
    $ws.Range("G3").Value = $epoch
    $ws.Range("H3").Value = $included

    $ws.Range("A4").Value = $file3
    $ws.Range("B4").Value = $readyHandoff
    $ws.Range("C4").Value = $xlf3
    $ws.Range("D4").Value = $t3
    $ws.Range("G4").Value = $epoch
    $ws.Range("H4").Value = $included

    $ws.Range("A5").Value = $file4
    $ws.Range("B5").Value = $readyHandoff
    $ws.Range("C5").Value = $xlf4
    $ws.Range("D5").Value = $t3
    $ws.Range("G5").Value = $epoch
    $ws.Range("H5").Value = $included

    $ws.Range("A6").Value = ".localization-config"
    $ws.Range("B6").Value = $notLocalized
    $ws.Range("D6").Value = $epoch
    $ws.Range("G6").Value = $epoch
    $ws.Range("H6").Value = $ignored

    $ws.Range("A2:A6").Style = "Hyperlink"
    $ws.Range("C2:C5").Style = "Hyperlink"

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), (MdUrl $file1), $null, $null, $file1) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), (XlfUrl $locale $xlf1), $null, $null, $xlf1) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), (MdUrl $file2), $null, $null, $file2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), (XlfUrl $locale $xlf2), $null, $null, $xlf2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A4"), (MdUrl $file3), $null, $null, $file3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C4"), (XlfUrl $locale $xlf3), $null, $null, $xlf3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A5"), (MdUrl $file4), $null, $null, $file4) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C5"), (XlfUrl $locale $xlf4), $null, $null, $xlf4) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A6"), $configUrl, $null, $null, ".localization-config") | Out-Null
}

$zh = $wb.Worksheets.Item("zh-cn")
FillDetailSheet $zh $xlf1zh $xlf2zh $xlf3zh $xlf4zh $zhTime1 $zhTime3 "zh"

$de = $wb.Worksheets.Item("de-de")
FillDetailSheet $de $xlf1de $xlf2de $xlf3de $xlf4de $deTime1 $deTime3 "de"
